# Data preprocessing pass over the "buku" (book) circulation category column:
# normalize/clean up the raw text values in column B (case, punctuation,
# whitespace) ahead of building the frequency matrix for CF input, and mark
# one of the "dirty" values with a hyperlink annotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order of these assignments matters for how the shared-string
# table gets rebuilt (new/changed strings are appended in the order they are
# first written), so keep this exact sequence: B2, B5, B4, B3, B8.
$ws.Range("B2").Value = "ekonomi, MiKro"
$ws.Range("B5").Value = "kepemimpinan; organisasi"
$ws.Range("B4").Value = "corel-draw12"
$ws.Range("B3").Value = "ekonomi!@##"
$ws.Range("B8").Value = "kepemimpinan, SEkolah, budaya     mutu"

# Flag the messy "ekonomi!@##" entry with a hyperlink annotation (adds the
# built-in Hyperlink cell style too).
$ws.Hyperlinks.Add($ws.Range("B3"), "https://example.com")

# Leave the selection where the user last clicked after finishing the edits.
$ws.Range("B16").Select()
